# Fruta / hortaliza, semanal
# Insert 7 new weekly data rows (2021-12-21, serial 44551) for Femacal de La Calera - Cereza
# above the existing row 255, shifting all subsequent rows down by 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows at 255..261 (pushes old 255.. down to 262..)
$ws.Rows("255:261").Insert()

# Common (unchanged across these new rows) column values
$commonA = 3
$commonB = "Femacal de La Calera"
$commonC = "Coquimbo"
$commonE = 5
$commonF = "Fruta"
$commonG = 100103
$commonH = "Frutos de hueso (carozo)"
$commonI = 100103001
$commonJ = "Cereza"
$commonQ = "`$/bandeja 10 kilos"
$commonR = "Provincia de Curicó"
$commonDate = 44551
$commonT = 10

function Set-CerezaRow {
    param($row, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg)

    $ws.Range("A$row").Value = $commonA
    $ws.Range("B$row").Value = $commonB
    $ws.Range("C$row").Value = $commonC
    $ws.Range("D$row").Value = $commonDate
    $ws.Range("E$row").Value = $commonE
    $ws.Range("F$row").Value = $commonF
    $ws.Range("G$row").Value = $commonG
    $ws.Range("H$row").Value = $commonH
    $ws.Range("I$row").Value = $commonI
    $ws.Range("J$row").Value = $commonJ
    $ws.Range("K$row").Value = $variedad
    $ws.Range("L$row").Value = $calidad
    $ws.Range("M$row").Value = $volumen
    $ws.Range("N$row").Value = $precioMin
    $ws.Range("O$row").Value = $precioMax
    $ws.Range("P$row").Value = $precioProm
    $ws.Range("Q$row").Value = $commonQ
    $ws.Range("R$row").Value = $commonR
    $ws.Range("S$row").Value = $precioKg
    $ws.Range("T$row").Value = $commonT
}

Set-CerezaRow 255 "Bing"    "Especial" 58 10000 10000 10000 1000
Set-CerezaRow 256 "Bing"    "Primera"  50 8000  8000  8000  800
Set-CerezaRow 257 "Brooks"  "Especial" 50 10000 10000 10000 1000
Set-CerezaRow 258 "Brooks"  "Primera"  54 8000  8000  8000  800
Set-CerezaRow 259 "Brooks"  "Segunda"  50 6000  6000  6000  600
Set-CerezaRow 260 "Rainier" "Especial" 50 12000 12000 12000 1200
Set-CerezaRow 261 "Rainier" "Primera"  48 10000 10000 10000 1000
